$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark column D as text so numeric-looking strings (e.g. "1.001")
# are stored as text, matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.191.37'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.848.69'
$ws.Range("D4").Value = '1.001'
$ws.Range("D5").Value = '313.08'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D7").Value = '0.4634'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.3696'
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").Value = '0.07274'
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("D10").Value = '0.8864'
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '0.07832'
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").Value = '1.877.95'
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("D14").Value = '5.390'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").Value = '6.508'
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").Value = '91.50'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '0.000008845'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("D20").Value = '27.219.36'
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("D22").Value = '5.058'
$ws.Range("D23").Value = '2.115.69'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").Value = '2.035'
$ws.Range("E25").Value = '  +9.00%  '
$ws.Range("D26").Value = '151.36'
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").Value = '18.41'
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").Value = '115.66'
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("D30").Value = '5.007'
$ws.Range("E30").Value = '  -2.61%  '
$ws.Range("D31").Value = '0.08835'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").Value = '3.131'
$ws.Range("E32").Value = '  +5.33%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '0.7809'
$ws.Range("E33").Value = '  +5.37%  '
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("D35").Value = '1.147'
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").Value = '2.689'
$ws.Range("E36").Value = '  +5.60%  '
$ws.Range("E37").Value = '  +2.06%  '
$ws.Range("D38").Value = '0.01942'
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("D39").Value = '0.05213'
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").Value = '2.953'
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("D41").Value = '7.029'
$ws.Range("D42").Value = '0.5040'
$ws.Range("E42").Value = '  -2.61%  '
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("D44").Value = '8.517'
$ws.Range("E44").Value = '  +2.91%  '
$ws.Range("D45").Value = '0.4759'
$ws.Range("E45").Value = '  -2.02%  '
$ws.Range("D46").Value = '10.40'
$ws.Range("E46").Value = '  +1.44%  '
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("D48").Value = '102.83'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").Value = '1.633'
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("D50").Value = '0.06191'
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").Value = '65.41'
$ws.Range("E51").Value = '  -0.01%  '

# Restore default cell style on column D so no stray number-format styling remains.
$ws.Range("D2:D51").Style = "Normal"
